$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.408.33'
$ws.Range("E2").Value = '  -3.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.693.17'
$ws.Range("E3").Value = '  -3.55%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.40'
$ws.Range("E5").Value = '  -2.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.21'
$ws.Range("E6").Value = '  -5.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.689.83'
$ws.Range("E7").Value = '  -3.48%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  -3.67%  '

$ws.Range("E11").Value = '  -3.95%  '

$ws.Range("E12").Value = '  -3.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.53'
$ws.Range("E13").Value = '  -6.04%  '

$ws.Range("E14").Value = '  -5.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.305.41'
$ws.Range("E15").Value = '  -3.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.684.82'
$ws.Range("E16").Value = '  -3.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.442.89'
$ws.Range("E17").Value = '  -3.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.57'
$ws.Range("E18").Value = '  +6.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.17'
$ws.Range("E19").Value = '  -3.78%  '

$ws.Range("E20").Value = '  -3.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '492.38'
$ws.Range("E21").Value = '  -2.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.17'
$ws.Range("E22").Value = '  -3.34%  '

$ws.Range("E23").Value = '  -1.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.59'
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("E25").Value = '  -5.88%  '

$ws.Range("E26").Value = '  -3.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.18'
$ws.Range("E27").Value = '  -3.43%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("E28").Value = '  -3.74%  '

$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("E30").Value = '  -1.64%  '

$ws.Range("E31").Value = '  -6.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.64'
$ws.Range("E32").Value = '  -3.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.54'
$ws.Range("E33").Value = '  -0.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.828.36'
$ws.Range("E34").Value = '  -3.56%  '

$ws.Range("E35").Value = '  -5.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.628.15'
$ws.Range("E36").Value = '  -3.51%  '

$ws.Range("E37").Value = '  +0.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.993'
$ws.Range("E38").Value = '  -4.45%  '

$ws.Range("E39").Value = '  -5.26%  '

$ws.Range("E40").Value = '  -7.02%  '

$ws.Range("E41").Value = '  -3.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '434.84'
$ws.Range("E42").Value = '  -10.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.52'
$ws.Range("E43").Value = '  -2.39%  '

$ws.Range("E44").Value = '  -5.95%  '

$ws.Range("E45").Value = '  -6.73%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.38'
$ws.Range("E46").Value = '  -1.63%  '

$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.51'
$ws.Range("E48").Value = '  -6.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '142.31'
$ws.Range("E49").Value = '  +1.77%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.752.72'
$ws.Range("E50").Value = '  -5.74%  '

$ws.Range("E51").Value = '  -3.80%  '
